$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1440975431801766
$ws.Range("C2").Value = 0.4700303102126987
$ws.Range("D2").Value = 0.5258244742912964
$ws.Range("E2").Value = 0.7251375554274488
$ws.Range("F2").Value = 0.7266482381404259

# Row 3
$ws.Range("B3").Value = 0.702818860278537
$ws.Range("C3").Value = 0.8472691904027576
$ws.Range("D3").Value = 4.825304880736022
$ws.Range("E3").Value = 2.196657661251753
$ws.Range("F3").Value = 2.127963838963247
$ws.Range("G3").Value = 23

# Row 4
$ws.Range("B4").Value = 0.2815506150436938
$ws.Range("C4").Value = 1.240068335160717
$ws.Range("D4").Value = 7.820836076261407
$ws.Range("E4").Value = 2.796575776956778
$ws.Range("F4").Value = 2.844899762820286
$ws.Range("G4").Value = 23

# Row 5
$ws.Range("B5").Value = 0.3158381654092163
$ws.Range("C5").Value = 1.158822516538064
$ws.Range("D5").Value = 7.724181108196982
$ws.Range("E5").Value = 2.779241102926657
$ws.Range("F5").Value = 2.823294572221369
$ws.Range("G5").Value = 23

# Row 6
$ws.Range("B6").Value = 0.3679265667964129
$ws.Range("C6").Value = 1.261171996547169
$ws.Range("D6").Value = 7.880840995669849
$ws.Range("E6").Value = 2.807283561678415
$ws.Range("F6").Value = 2.8456173148262
$ws.Range("G6").Value = 23

# Row 7
$ws.Range("B7").Value = 0.3339394964461404
$ws.Range("C7").Value = 1.329628849026482
$ws.Range("D7").Value = 8.057429588549642
$ws.Range("E7").Value = 2.838561182808932
$ws.Range("F7").Value = 2.882202614392551
$ws.Range("G7").Value = 23

# Row 8
$ws.Range("B8").Value = 0.2421922022797161
$ws.Range("C8").Value = 1.326115933640845
$ws.Range("D8").Value = 8.153426962725963
$ws.Range("E8").Value = 2.85542062798565
$ws.Range("F8").Value = 2.909074420883571
$ws.Range("G8").Value = 23

# Row 9
$ws.Range("B9").Value = 0.2908000157509625
$ws.Range("C9").Value = 1.361567814352823
$ws.Range("D9").Value = 8.088431141923708
$ws.Range("E9").Value = 2.844016726730648
$ws.Range("F9").Value = 2.892694004914857
$ws.Range("G9").Value = 23

# Row 10
$ws.Range("B10").Value = 0.2584427042972859
$ws.Range("C10").Value = 1.344876645801738
$ws.Range("D10").Value = 8.092074104027965
$ws.Range("E10").Value = 2.844657115370492
$ws.Range("F10").Value = 2.896561236035992
$ws.Range("G10").Value = 23

# Row 11
$ws.Range("B11").Value = 0.1521280651501415
$ws.Range("C11").Value = 1.264221546682911
$ws.Range("D11").Value = 7.809337215816921
$ws.Range("E11").Value = 2.794519138566942
$ws.Range("F11").Value = 2.853088184558203
$ws.Range("G11").Value = 23
